# 自动更新Excel文件 - daily decrement of remaining days (剩余, column E)
# Rule: each row's "剩余" (remaining days, col E) counts down by 1 per day.
# When a row reaches 1 (i.e. the last remaining day), it is treated as
# restocked: 剩余 resets back to the row's "总天" (total days, col D) and the
# "开始时间" (start date, col F, stored as a plain YYYYMMDD integer) is pushed
# forward by 10 days.
# Row 36 is left untouched (its data is already inconsistent - E36 == D36 -
# and is never part of the automatic update).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 2 }

for ($row = 2; $row -le $lastRow; $row++) {

    if ($row -eq 36) { continue }

    $totalDays = $ws.Cells.Item($row, 4).Value2
    $remaining = $ws.Cells.Item($row, 5).Value2
    $startDate = $ws.Cells.Item($row, 6).Value2

    if ($remaining -eq $null -or $totalDays -eq $null) { continue }

    if ($remaining -eq 1) {
        # Restock: reset remaining to total days, push start date 10 days out.
        $ws.Cells.Item($row, 5).Value = $totalDays

        $year = [Math]::Floor($startDate / 10000)
        $month = [Math]::Floor(($startDate % 10000) / 100)
        $day = $startDate % 100

        $dt = (Get-Date -Year $year -Month $month -Day $day).AddDays(10)
        $newStartDate = ($dt.Year * 10000) + ($dt.Month * 100) + $dt.Day

        $ws.Cells.Item($row, 6).Value = $newStartDate
    } else {
        $ws.Cells.Item($row, 5).Value = $remaining - 1
    }
}
